$wb = $excel.ActiveWorkbook

# Rename the existing (first) sheet to "Testing"
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Testing"

# Add the text value to F7 of the Testing sheet
$ws.Range("F7").Value = "This is a test for check and control versions in excel with GitHub :)"

# Add a new worksheet after the existing one, named "newSheet"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "newSheet"

# Return focus to the Testing sheet, with the selection having moved to F8
# (as if the user had just finished typing into F7)
$ws.Activate() | Out-Null
$ws.Range("F8").Select() | Out-Null
